# Lays results out as a table.
#
# Source paragraphs (before):
#   1: "Weather:  Sunny"
#   2: "Location: Edmonton"
#   3: "Date: March 1"
#
# Target paragraphs (after):
#   1: "Authour: Joe"                                   (single run)
#   2: "Location: " + "New York"                        (two runs)
#   3: "Date: " + "August" + " 1"                        (three runs)
#   4: "End"                                             (new paragraph, single run)

$d = $word.ActiveDocument

# --- Paragraph 1: plain text swap, single run stays a single run. ---
$d.Content.Find.Execute("Weather:  Sunny", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Authour: Joe", 2)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphRuns($paragraph, $runsXml) {
    # Replace a paragraph's content with an explicit run layout (keeps
    # each supplied <w:r> as its own run instead of Word's usual
    # same-format run coalescing) by round-tripping through WordOpenXML.
    $full = $paragraph.Range
    $sub = $d.Range($full.Start, $full.End - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $ns + '><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $sub.InsertXML($pkg)
}

# --- Paragraph 2: "Location: Edmonton" -> "Location: " + "New York" ---
$p2 = $d.Paragraphs(2)
$p2Runs = '<w:r><w:t xml:space="preserve">Location: </w:t></w:r>' +
          '<w:r><w:t>New York</w:t></w:r>'
Set-ParagraphRuns $p2 $p2Runs

# --- Paragraph 3: "Date: March 1" -> "Date: " + "August" + " 1" ---
$p3 = $d.Paragraphs(3)
$p3Runs = '<w:r><w:t xml:space="preserve">Date: </w:t></w:r>' +
          '<w:r><w:t>August</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> 1</w:t></w:r>'
Set-ParagraphRuns $p3 $p3Runs

# --- New paragraph 4: "End" ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newFull = $newPara.Range
$newSub = $d.Range($newFull.Start, $newFull.End - 1)
$newSub.Text = "End"

Write-Host "Final paragraphs:"
foreach ($p in $d.Paragraphs) {
    Write-Host ("  [" + $p.Range.Text + "]")
}
